$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "_x597D_"
$ws.Range("A16").Value = "_x597d_"
$ws.Range("A17").Value = "_x597G_"
$ws.Range("A18").Value = "_x_x_x"
